$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.308.74'

$ws.Range('E2').Value = '  -0.11%  '

$ws.Range('D3').Value = '1.931.31'

$ws.Range('E3').Value = '  -0.11%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7568'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').Value = '  +5.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.71'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').Value = '  -2.50%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3177'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').Value = '  -2.84%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.53'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').Value = '  -0.02%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07005'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').Value = '  -2.46%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7800'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').Value = '  -2.38%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07990'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').Value = '  -1.13%  '

$ws.Range('D13').Value = '1.930.28'

$ws.Range('E13').Value = '  -0.09%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.362'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').Value = '  -0.99%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.25'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').Value = '  -0.27%  '

$ws.Range('E16').Value = '  -2.37%  '

$ws.Range('D17').Value = '30.309.14'

$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '252.62'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').Value = '  +0.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007917'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').Value = '  -2.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.738'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').Value = '  -0.81%  '

$ws.Range('D21').Value = '2.182.18'

$ws.Range('E21').Value = '  -0.06%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').Value = '  -0.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').Value = '  -0.18%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.666'
$ws.Range('D24').ClearFormats()

$ws.Range('E24').Value = '  -3.54%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.497'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').Value = '  -2.42%  '

$ws.Range('E26').Value = '  +0.25%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1337'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').Value = '  +4.08%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.94'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').Value = '  -1.50%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.207'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').Value = '  -5.20%  '

$ws.Range('E30').Value = '  +0.24%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.513'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').Value = '  -1.87%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.375'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').Value = '  -1.06%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.116'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').Value = '  -2.10%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05159'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').Value = '  -0.77%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.287'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').Value = '  +1.73%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7460'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').Value = '  -0.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.772'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').Value = '  +0.07%  '

$ws.Range('E38').Value = '  -0.71%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '77.53'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').Value = '  -1.72%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.404'
$ws.Range('D41').ClearFormats()

$ws.Range('E41').Value = '  -0.72%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4463'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').Value = '  -1.43%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.965'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').Value = '  -2.96%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').Value = '  -0.08%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8309'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').Value = '  -1.10%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.79'
$ws.Range('D46').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.732'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').Value = '  -0.35%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.465'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').Value = '  +0.87%  '

$ws.Range('B49').Value = 'Elrond'

$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.37'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').Value = '  +2.07%  '

$ws.Range('B50').Value = 'Maker'

$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '981.25'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').Value = '  +11.12%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06006'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').Value = '  -0.98%  '
